$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.012.18'
$ws.Range("E2").Value = '  -3.15%  '

$ws.Range("D3").Value = '3.162.55'
$ws.Range("E3").Value = '  -8.49%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.86'
$ws.Range("E5").Value = '  -4.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.93'
$ws.Range("E6").Value = '  -2.58%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '3.157.60'
$ws.Range("E9").Value = '  -8.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  -6.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.58'
$ws.Range("E11").Value = '  -4.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  -4.17%  '

$ws.Range("D13").Value = '3.715.15'
$ws.Range("E13").Value = '  -8.39%  '

$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.38'
$ws.Range("E15").Value = '  -5.88%  '

$ws.Range("D16").Value = '64.168.68'
$ws.Range("E16").Value = '  -2.93%  '

$ws.Range("E17").Value = '  -5.89%  '

$ws.Range("D18").Value = '3.174.98'
$ws.Range("E18").Value = '  -7.87%  '

$ws.Range("E19").Value = '  -5.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.93'
$ws.Range("E20").Value = '  -6.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.18'
$ws.Range("E21").Value = '  -4.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.17'
$ws.Range("E22").Value = '  -6.77%  '

$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.83'
$ws.Range("E24").Value = '  -5.39%  '

$ws.Range("E25").Value = '  -7.01%  '

$ws.Range("E26").Value = '  -4.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.36'
$ws.Range("E27").Value = '  -4.35%  '

$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.56'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.89'
$ws.Range("E32").Value = '  -4.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.03'
$ws.Range("E33").Value = '  -7.53%  '

$ws.Range("E34").Value = '  -6.66%  '

$ws.Range("E35").Value = '  -9.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.69'
$ws.Range("E36").Value = '  -2.11%  '

$ws.Range("E37").Value = '  -7.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.803'
$ws.Range("E38").Value = '  -8.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.88'
$ws.Range("E39").Value = '  -10.69%  '

$ws.Range("E40").Value = '  -5.81%  '

$ws.Range("E41").Value = '  -7.19%  '

$ws.Range("D42").Value = '2.635.57'
$ws.Range("E42").Value = '  -4.62%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.12'
$ws.Range("E43").Value = '  -7.76%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.95'
$ws.Range("E44").Value = '  -7.65%  '

$ws.Range("E45").Value = '  -4.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.46'
$ws.Range("E46").Value = '  -4.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.35'
$ws.Range("E47").Value = '  -1.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.28'
$ws.Range("E48").Value = '  -4.47%  '

$ws.Range("E49").Value = '  -7.34%  '

$ws.Range("E50").Value = '  -0.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.05%  '
